$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.385.48'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '2.636.72'
$ws.Range('E3').Value = '  +1.30%  '
$ws.Range('E4').Value = '  +0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '536.69'
$c.Style = 'Normal'

$ws.Range('E5').Value = '  -0.77%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '145.09'
$c.Style = 'Normal'

$ws.Range('E6').Value = '  +2.51%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +1.29%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '7.05'
$c.Style = 'Normal'

$ws.Range('E9').Value = '  +9.55%  '
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D13').Value = '3.101.73'
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('D14').Value = '59.301.70'
$ws.Range('E14').Value = '  +0.23%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '21.25'
$c.Style = 'Normal'

$ws.Range('E15').Value = '  +3.29%  '
$ws.Range('D16').Value = '2.655.04'
$ws.Range('E16').Value = '  +2.14%  '
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('E18').Value = '  +3.36%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '338.41'
$c.Style = 'Normal'

$ws.Range('E19').Value = '  -0.82%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '10.28'
$c.Style = 'Normal'

$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('E21').Value = '  -2.38%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '66.23'
$c.Style = 'Normal'

$ws.Range('E23').Value = '  -2.06%  '
$ws.Range('E24').Value = '  +1.92%  '
$ws.Range('E25').Value = '  -0.05%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.01'
$c.Style = 'Normal'

$ws.Range('E26').Value = '  +0.75%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '7.30'
$c.Style = 'Normal'

$ws.Range('E27').Value = '  +1.50%  '
$ws.Range('E28').Value = '  +0.17%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'

$ws.Range('E29').Value = '  -0.06%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.66'
$c.Style = 'Normal'

$ws.Range('E30').Value = '  -1.80%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '5.89'
$c.Style = 'Normal'

$ws.Range('E31').Value = '  +1.32%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '151.50'
$c.Style = 'Normal'

$ws.Range('E32').Value = '  +1.22%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '18.83'
$c.Style = 'Normal'

$ws.Range('E33').Value = '  +0.57%  '
$ws.Range('E35').Value = '  +2.20%  '
$ws.Range('E36').Value = '  +1.97%  '
$ws.Range('E37').Value = '  +0.69%  '
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('E39').Value = '  +1.19%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '284.87'
$c.Style = 'Normal'

$ws.Range('E40').Value = '  +3.57%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('E44').Value = '  +2.97%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '19.11'
$c.Style = 'Normal'

$ws.Range('E45').Value = '  +2.60%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.0942'
$c.Style = 'Normal'

$ws.Range('E46').Value = '  -1.34%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0226'
$c.Style = 'Normal'

$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('D48').Value = '1.959.83'
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '18.37'
$c.Style = 'Normal'

$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '4.55'
$c.Style = 'Normal'

$ws.Range('E50').Value = '  +0.53%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '111.48'
$c.Style = 'Normal'

$ws.Range('E51').Value = '  +0.04%  '
